$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 1.3
$ws.Range("H7").Value = 5.5
$ws.Range("I7").Value = 9
$ws.Range("K7").Value = 2.38
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.85
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.75
$ws.Range("W7").Value = 5.5
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 7.5
$ws.Range("AJ7").Value = 29
$ws.Range("AR7").Value = 51
$ws.Range("AT7").Value = 2.75
$ws.Range("AW7").Value = 10
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 3.2
$ws.Range("L8").Value = 4
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("X8").Value = 10
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 21
$ws.Range("AD8").Value = 6.5
$ws.Range("AI8").Value = 15
$ws.Range("AK8").Value = 34
$ws.Range("AW8").Value = 5
$ws.Range("AY8").Value = 34
$ws.Range("AZ8").Value = 67
$ws.Range("BB8").Value = 301
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 2.88
$ws.Range("J9").Value = 3.2
$ws.Range("O9").Value = 1.53
$ws.Range("P9").Value = 2.38
$ws.Range("Q9").Value = 2.7
$ws.Range("R9").Value = 1.44
$ws.Range("U9").Value = 2.2
$ws.Range("V9").Value = 1.62
$ws.Range("AA9").Value = 23
$ws.Range("AC9").Value = 6
$ws.Range("AD9").Value = 5.5
$ws.Range("AM9").Value = 51
$ws.Range("AO9").Value = 15
$ws.Range("AS9").Value = 301
$ws.Range("AU9").Value = 9.5
$ws.Range("AV9").Value = 81
$ws.Range("AZ9").Value = 81
$ws.Range("BA9").Value = 126
$ws.Range("G10").Value = 1.42
$ws.Range("H10").Value = 4.2
$ws.Range("I10").Value = 8
$ws.Range("J10").Value = 2
$ws.Range("L10").Value = 8
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("W10").Value = 5.5
$ws.Range("AE10").Value = 26
$ws.Range("AH10").Value = 13
$ws.Range("AJ10").Value = 23
$ws.Range("AM10").Value = 67
$ws.Range("AQ10").Value = 23
$ws.Range("AS10").Value = 201
$ws.Range("AZ10").Value = 201
$ws.Range("BA10").Value = 251
$ws.Range("AQ11").Value = 29
$ws.Range("BB11").Value = 351
$ws.Range("J12").Value = 4
$ws.Range("O12").Value = 1.23
$ws.Range("P12").Value = 3.4
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 1.93
$ws.Range("V12").Value = 2.05
$ws.Range("W12").Value = 13
$ws.Range("AB12").Value = 32
$ws.Range("AC12").Value = 11.5
$ws.Range("AD12").Value = 6.8
$ws.Range("AE12").Value = 13
$ws.Range("AG12").Value = 350
$ws.Range("AH12").Value = 8
$ws.Range("AI12").Value = 9.5
$ws.Range("AK12").Value = 16.5
$ws.Range("AL12").Value = 14.5
$ws.Range("AM12").Value = 23
$ws.Range("AN12").Value = 5.6
$ws.Range("AQ12").Value = 100
$ws.Range("AR12").Value = 120
$ws.Range("AT12").Value = 2.85
$ws.Range("AU12").Value = 6.7
$ws.Range("AV12").Value = 55
$ws.Range("AW12").Value = 3.85
$ws.Range("AY12").Value = 17
$ws.Range("BB12").Value = 200
$ws.Range("G19").Value = 1.62
$ws.Range("H19").Value = 3.55
$ws.Range("I19").Value = 5.3
$ws.Range("J19").Value = 2.15
$ws.Range("K19").Value = 2.15
$ws.Range("L19").Value = 5.3
$ws.Range("M19").Value = 1.01
$ws.Range("N19").Value = 9.15
$ws.Range("O19").Value = 1.29
$ws.Range("P19").Value = 3
$ws.Range("Q19").Value = 1.91
$ws.Range("R19").Value = 1.8
$ws.Range("S19").Value = 1.38
$ws.Range("T19").Value = 2.6
$ws.Range("U19").Value = 1.85
$ws.Range("V19").Value = 1.75
$ws.Range("X19").Value = 7.2
$ws.Range("Z19").Value = 11.75
$ws.Range("AA19").Value = 13.5
$ws.Range("AC19").Value = 9.5
$ws.Range("AD19").Value = 7
$ws.Range("AE19").Value = 16.5
$ws.Range("AF19").Value = 80
$ws.Range("AG19").Value = 700
$ws.Range("AH19").Value = 13.5
$ws.Range("AI19").Value = 32
$ws.Range("AJ19").Value = 17
$ws.Range("AK19").Value = 110
$ws.Range("AL19").Value = 55
$ws.Range("AM19").Value = 55
$ws.Range("AN19").Value = 3.4
$ws.Range("AO19").Value = 7.7
$ws.Range("AP19").Value = 17
$ws.Range("AQ19").Value = 25
$ws.Range("AR19").Value = 55
$ws.Range("AS19").Value = 250
$ws.Range("AT19").Value = 2.6
$ws.Range("AU19").Value = 7.4
$ws.Range("AV19").Value = 70
$ws.Range("AW19").Value = 6.8
$ws.Range("AX19").Value = 30
$ws.Range("AY19").Value = 32
$ws.Range("AZ19").Value = 175
$ws.Range("BA19").Value = 200
$ws.Range("BB19").Value = 400
